# "annotated the 27 landmarks"
# Fill in the previously-empty point_id values in column B (and correct a
# few that had placeholder/incorrect numbers) for the landmark_points sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    5  = 138
    6  = 5752
    7  = 7395
    8  = 4503
    9  = 8614
    10 = 9396
    11 = 4328
    12 = 8232
    13 = 6703
    14 = 3171
    15 = 5549
    16 = 2764
    17 = 10595
    18 = 6576
    19 = 6274
    24 = 1706
    25 = 281
    26 = 6280
    27 = 9124
    28 = 6083
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}

# Move the selection to where the author's cursor ended up after the edits.
[void]$ws.Range("F6").Select()
